$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 3 ("CustomShape 3"): "sceneFormat" -> "contentFormat"
$shp3 = $s.Shapes.Item(3)
$tr3 = $shp3.TextFrame.TextRange
$tr3.Characters(1, 11).Text = "contentFormat"

# Shape 4 ("CustomShape 4"): "sceneLength" -> two runs "content" + "Length"
$shp4 = $s.Shapes.Item(4)
$tr4 = $shp4.TextFrame.TextRange
$tr4.Characters(1, 11).Text = ""
$para4 = $tr4.Paragraphs(1)
$run4a = $para4.InsertAfter("content")
$run4a.LanguageID = "en-US"
$run4a.Font.Size = 12
$run4a.Font.Name = "Courier New"
$run4a.Font.Color.RGB = 0
$run4b = $run4a.InsertAfter("Length")
$run4b.LanguageID = "en-US"
$run4b.Font.Size = 12
$run4b.Font.Name = "Courier New"
$run4b.Font.Color.RGB = 0

# Shape 14 ("CustomShape 14"): "scene" -> "content"
$shp14 = $s.Shapes.Item(14)
$tr14 = $shp14.TextFrame.TextRange
$tr14.Characters(1, 5).Text = "content"
